# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates columns I (DAMSLTag) and J (DialogAct) for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=6;   I='sd'; J='Statement-non-opinion'},
    @{Row=7;   I='sd'; J='Statement-non-opinion'},
    @{Row=15;  I='aa'; J='Agree/Accept'},
    @{Row=25;  I='sv'; J='Statement-opinion'},
    @{Row=42;  I='sd'; J='Statement-non-opinion'},
    @{Row=53;  I='%';  J='Uninterpretable'},
    @{Row=54;  I='sv'; J='Statement-opinion'},
    @{Row=55;  I='aa'; J='Agree/Accept'},
    @{Row=58;  I='sd'; J='Statement-non-opinion'},
    @{Row=64;  I='sv'; J='Statement-opinion'},
    @{Row=95;  I='sv'; J='Statement-opinion'},
    @{Row=105; I='sd'; J='Statement-non-opinion'},
    @{Row=106; I='sd'; J='Statement-non-opinion'},
    @{Row=111; I='aa'; J='Agree/Accept'},
    @{Row=152; I='aa'; J='Agree/Accept'},
    @{Row=158; I='sd'; J='Statement-non-opinion'},
    @{Row=159; I='sd'; J='Statement-non-opinion'},
    @{Row=172; I='sv'; J='Statement-opinion'},
    @{Row=174; I='sv'; J='Statement-opinion'},
    @{Row=177; I='sv'; J='Statement-opinion'},
    @{Row=185; I='sd'; J='Statement-non-opinion'},
    @{Row=188; I='aa'; J='Agree/Accept'},
    @{Row=190; I='sd'; J='Statement-non-opinion'},
    @{Row=194; I='aa'; J='Agree/Accept'},
    @{Row=200; I='aa'; J='Agree/Accept'},
    @{Row=202; I='aa'; J='Agree/Accept'},
    @{Row=231; I='sd'; J='Statement-non-opinion'},
    @{Row=237; I='sd'; J='Statement-non-opinion'},
    @{Row=244; I='ba'; J='Appreciation'},
    @{Row=249; I='ba'; J='Appreciation'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I   # Column I = DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.J  # Column J = DialogAct
}
